# Weekly Fruta/Hortaliza update: insert a new daily price record as row 285
# (shifting the existing rows 285:363 down to 286:364) on the "Cilantro -
# Vega Modelo de Temuco" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 285..363 down by one to make room for the new record.
$ws.Rows("285:285").Insert()

# Populate the newly inserted row 285 with the new observation.
$ws.Cells.Item(285, 1).Value = 10
$ws.Cells.Item(285, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(285, 3).Value = "La Araucanía"
$ws.Cells.Item(285, 4).Value = 44736
$ws.Cells.Item(285, 5).Value = 9
$ws.Cells.Item(285, 6).Value = 100112040
$ws.Cells.Item(285, 7).Value = "Cilantro"
$ws.Cells.Item(285, 8).Value = "Sin especificar"
$ws.Cells.Item(285, 9).Value = "Primera"
$ws.Cells.Item(285, 10).Value = 40
$ws.Cells.Item(285, 11).Value = 5000
$ws.Cells.Item(285, 12).Value = 5000
$ws.Cells.Item(285, 13).Value = 5000
$ws.Cells.Item(285, 14).Value = "$/docena de atados (2 kilos)"
$ws.Cells.Item(285, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(285, 16).Value = 2500
$ws.Cells.Item(285, 17).Value = 2
$ws.Cells.Item(285, 18).Value = "Hortaliza"
